$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("B2").Value = -1
$ws.Range("C2").Value = "date"
$ws.Range("E2").Value = 209
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = -1
$ws.Range("L2").Value = -1

# --- Row 3 ---
$ws.Range("B3").Value = -1
$ws.Range("C3").Value = "date"

# --- Row 4 ---
$ws.Range("B4").Value = -1
$ws.Range("C4").Value = "date"

# --- Row 5 ---
$ws.Range("B5").Value = -1
$ws.Range("C5").Value = "date"

# --- Row 6 ---
$ws.Range("B6").Value = -1
$ws.Range("C6").Value = "date"

# --- Row 7 ---
$ws.Range("B7").Value = -1
$ws.Range("C7").Value = "date"
$ws.Range("U7").Value = 0

# --- Row 8 ---
$ws.Range("B8").Value = -1
$ws.Range("C8").Value = "date"
$ws.Range("U8").Value = 0

# --- Row 9 ---
$ws.Range("B9").Value = -1
$ws.Range("C9").Value = "date"
$ws.Range("U9").Value = 0

# --- Row 10 ---
$ws.Range("B10").Value = -1
$ws.Range("C10").Value = "date"

# --- Row 11 ---
$ws.Range("B11").Value = -1
$ws.Range("C11").Value = "date"

# --- Row 12 ---
$ws.Range("B12").Value = -1
$ws.Range("C12").Value = "date"

# --- Row 13 ---
$ws.Range("B13").Value = -1
$ws.Range("C13").Value = "integer"
$ws.Range("E13").Value = 0

# --- Row 14 (new row) ---
# A14 and S14 are empty-text cells (not blank cells). A lone apostrophe is
# Excel's "treat as text" quote-prefix marker and collapses to an empty
# string value, which is exactly what's needed here; resetting the style
# afterwards drops the quote-prefix formatting flag it leaves behind.
$ws.Range("A14").Value = "'"
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").Value = -1
$ws.Range("C14").Value = "date"
$ws.Range("D14").Value = 209
$ws.Range("E14").Value = 209
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = -1
$ws.Range("L14").Value = -1
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 19
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 207
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = "<Unspecified>"
$ws.Range("S14").Value = "'"
$ws.Range("S14").Style = "Normal"
$ws.Range("T14").Value = 207
$ws.Range("U14").Value = 0
$ws.Range("V14").Value = 38
$ws.Range("W14").Value = 1
$ws.Range("X14").Value = "['not on junior sheet']"
